$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix apartment name typo: add a space ---
$ws.Range("B59").Value = "해링턴플레이스 감삼3차"

# --- Data corrections across many rows (E/F totals + per-month buckets) ---
$ws.Range("E9").Value = 268
$ws.Range("F9").Value = 0.99260000000000004
$ws.Range("AE9").Value = 3
$ws.Range("E21").Value = 1852
$ws.Range("F21").Value = 0.99139999999999995
$ws.Range("AD21").Value = 2
$ws.Range("E23").Value = 217
$ws.Range("F23").Value = 0.95179999999999998
$ws.Range("R23").Value = 1
$ws.Range("E24").Value = 273
$ws.Range("F24").Value = 0.45650000000000002
$ws.Range("I24").Value = 83
$ws.Range("E25").Value = 454
$ws.Range("F25").Value = 0.88670000000000004
$ws.Range("W25").Value = 1
$ws.Range("E30").Value = 223
$ws.Range("F30").Value = 0.95299999999999996
$ws.Range("AK30").Value = 1
$ws.Range("E35").Value = 909
$ws.Range("F35").Value = 0.99560000000000004
$ws.Range("AD35").Value = 1
$ws.Range("E54").Value = 114
$ws.Range("F54").Value = 0.13619999999999999
$ws.Range("M54").Value = 52
$ws.Range("N54").Value = 9
$ws.Range("E55").Value = 655
$ws.Range("F55").Value = 0.97330000000000005
$ws.Range("R55").Value = 1
$ws.Range("E56").Value = 403
$ws.Range("F56").Value = 0.88380000000000003
$ws.Range("Q56").Value = 1
$ws.Range("E57").Value = 611
$ws.Range("F57").Value = 0.76949999999999996
$ws.Range("U57").Value = 7
$ws.Range("E58").Value = 182
$ws.Range("F58").Value = 0.2984
$ws.Range("L58").Value = 36
$ws.Range("E60").Value = 491
$ws.Range("F60").Value = 0.94420000000000004
$ws.Range("P60").Value = 1
$ws.Range("E66").Value = 735
$ws.Range("F66").Value = 0.98660000000000003
$ws.Range("AD66").Value = 1
$ws.Range("E70").Value = 249
$ws.Range("F70").Value = 0.96509999999999996
$ws.Range("AQ70").Value = 1
$ws.Range("E71").Value = 712
$ws.Range("F71").Value = 0.99719999999999998
$ws.Range("AA71").Value = 1
$ws.Range("E86").Value = 116
$ws.Range("F86").Value = 0.9667
$ws.Range("AG86").Value = 1
$ws.Range("E88").Value = 907
$ws.Range("F88").Value = 0.98160000000000003
$ws.Range("S88").Value = 3
$ws.Range("E89").Value = 599
$ws.Range("F89").Value = 0.89539999999999997
$ws.Range("S89").Value = 5
$ws.Range("E90").Value = 312
$ws.Range("F90").Value = 0.41320000000000001
$ws.Range("Q90").Value = 5
$ws.Range("E93").Value = 312
$ws.Range("F93").Value = 0.7429
$ws.Range("T93").Value = 1
$ws.Range("E94").Value = 59
$ws.Range("F94").Value = 0.2576
$ws.Range("Q94").Value = 9
$ws.Range("E95").Value = 506
$ws.Range("F95").Value = 0.87849999999999995
$ws.Range("Q95").Value = 2
$ws.Range("E101").Value = 251
$ws.Range("F101").Value = 0.95440000000000003
$ws.Range("AF101").Value = 1
$ws.Range("E104").Value = 59
$ws.Range("F104").Value = 0.89390000000000003
$ws.Range("AV104").Value = 1
$ws.Range("E109").Value = 406
$ws.Range("F109").Value = 0.86570000000000003
$ws.Range("W109").Value = 2
$ws.Range("E111").Value = 999
$ws.Range("F111").Value = 0.91739999999999999
$ws.Range("Q111").Value = 5
$ws.Range("E112").Value = 456
$ws.Range("F112").Value = 0.45190000000000002
$ws.Range("K112").Value = 32
$ws.Range("E113").Value = 21
$ws.Range("F113").Value = 0.42
$ws.Range("N113").Value = 21
$ws.Range("E114").Value = 303
$ws.Range("F114").Value = 0.83009999999999995
$ws.Range("P114").Value = 5
$ws.Range("E133").Value = 14
$ws.Range("F133").Value = 0.058299999999999998
$ws.Range("AC133").Value = 1
$ws.Range("E139").Value = 390
$ws.Range("F139").Value = 0.99490000000000001
$ws.Range("Y139").Value = 1
$ws.Range("E140").Value = 1193
$ws.Range("F140").Value = 0.98839999999999995
$ws.Range("Z140").Value = 3
$ws.Range("E141").Value = 174
$ws.Range("F141").Value = 0.90629999999999999
$ws.Range("AP141").Value = 1
$ws.Range("E145").Value = 208
$ws.Range("F145").Value = 0.80620000000000003
$ws.Range("P145").Value = 3
$ws.Range("E147").Value = 416
$ws.Range("F147").Value = 0.82379999999999998
$ws.Range("O147").Value = 10
$ws.Range("E150").Value = 28
$ws.Range("F150").Value = 0.71789999999999998
$ws.Range("U150").Value = 2
$ws.Range("E152").Value = 646
$ws.Range("F152").Value = 0.89970000000000006
$ws.Range("Q152").Value = 8
$ws.Range("E153").Value = 3
$ws.Range("F153").Value = 0.024400000000000002
$ws.Range("R153").Value = 2
$ws.Range("E156").Value = 950
$ws.Range("F156").Value = 0.99060000000000004
$ws.Range("AB156").Value = 1
$ws.Range("E179").Value = 535
$ws.Range("F179").Value = 0.96750000000000003
$ws.Range("AT179").Value = 1
$ws.Range("E182").Value = 1084
$ws.Range("F182").Value = 0.94259999999999999
$ws.Range("X182").Value = 1
$ws.Range("E184").Value = 494
$ws.Range("F184").Value = 0.91310000000000002
$ws.Range("T184").Value = 4
$ws.Range("E185").Value = 12
$ws.Range("F185").Value = 0.034200000000000001
$ws.Range("M185").Value = 3
$ws.Range("E188").Value = 336
$ws.Range("F188").Value = 0.89359999999999995
$ws.Range("AD188").Value = 1
$ws.Range("E197").Value = 60
$ws.Range("F197").Value = 0.35289999999999999
$ws.Range("AD197").Value = 2
$ws.Range("E205").Value = 191
$ws.Range("F205").Value = 0.77639999999999998
$ws.Range("Y205").Value = 3
$ws.Range("E213").Value = 84
$ws.Range("F213").Value = 0.71789999999999998
$ws.Range("L213").Value = 5
$ws.Range("E216").Value = 116
$ws.Range("F216").Value = 0.2944
$ws.Range("Z216").Value = 10
$ws.Range("E221").Value = 614
$ws.Range("F221").Value = 0.98080000000000001
$ws.Range("AN221").Value = 2
$ws.Range("E222").Value = 538
$ws.Range("F222").Value = 0.94720000000000004
$ws.Range("AD222").Value = 1
$ws.Range("E225").Value = 294
$ws.Range("F225").Value = 0.91879999999999995
$ws.Range("AL225").Value = 1
$ws.Range("E234").Value = 81
$ws.Range("F234").Value = 0.9
$ws.Range("AN234").Value = 2
$ws.Range("E235").Value = 530
$ws.Range("F235").Value = 0.93310000000000004
$ws.Range("Z235").Value = 7
$ws.Range("E238").Value = 298
$ws.Range("F238").Value = 0.215
$ws.Range("I238").Value = 270
$ws.Range("E239").Value = 401
$ws.Range("F239").Value = 0.88329999999999997
$ws.Range("V239").Value = 1
$ws.Range("E242").Value = 387
$ws.Range("F242").Value = 0.89380000000000004
$ws.Range("S242").Value = 1
$ws.Range("E243").Value = 85
$ws.Range("F243").Value = 0.46700000000000003
$ws.Range("Q243").Value = 5
$ws.Range("E245").Value = 381
$ws.Range("F245").Value = 0.37909999999999999
$ws.Range("I245").Value = 116
$ws.Range("E246").Value = 102
$ws.Range("F246").Value = 0.47220000000000001
$ws.Range("P246").Value = 2
$ws.Range("E247").Value = 116089
$ws.Range("F247").Value = 0.90539999999999998
$ws.Range("I247").Value = 24586
$ws.Range("K247").Value = 17243
$ws.Range("L247").Value = 9440
$ws.Range("M247").Value = 5475
$ws.Range("N247").Value = 3547
$ws.Range("O247").Value = 3142
$ws.Range("P247").Value = 2611
$ws.Range("Q247").Value = 1151
$ws.Range("R247").Value = 1020
$ws.Range("S247").Value = 910
$ws.Range("T247").Value = 1004
$ws.Range("U247").Value = 791
$ws.Range("V247").Value = 411
$ws.Range("W247").Value = 349
$ws.Range("X247").Value = 311
$ws.Range("Y247").Value = 226
$ws.Range("Z247").Value = 508
$ws.Range("AA247").Value = 213
$ws.Range("AB247").Value = 201
$ws.Range("AC247").Value = 143
$ws.Range("AD247").Value = 150
$ws.Range("AE247").Value = 73
$ws.Range("AF247").Value = 58
$ws.Range("AG247").Value = 43
$ws.Range("AK247").Value = 29
$ws.Range("AL247").Value = 21
$ws.Range("AN247").Value = 35
$ws.Range("AP247").Value = 16
$ws.Range("AQ247").Value = 14
$ws.Range("AT247").Value = 9
$ws.Range("AV247").Value = 8

# --- Row 247 (T247) crosses 1000 -> apply thousands-separator format like its neighbours ---
$ws.Range("T247").NumberFormat = "#,##0"

# --- Restore the active-cell selection shown in the sheet view ---
$ws.Range("H13").Select() | Out-Null
